$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values per the diff:
# Rows 2-25  -> 7318
# Rows 26-53 -> 7310
# Rows 104-145 -> 7310

$ws.Range("C2:C25").Value = 7318
$ws.Range("C26:C53").Value = 7310
$ws.Range("C104:C145").Value = 7310
